# Scheduled market-data refresh: update Leve profit-calc columns (H-N) for the
# rows whose underlying item prices changed. Generated from the upstream diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 11327.513
$ws.Range("I62").Value = 12978.046
$ws.Range("J62").Value = 9191.529
$ws.Range("K62").Value = 12978.046
$ws.Range("L62").Value = 9191.529
$ws.Range("M62").Value = -12354.046
$ws.Range("N62").Value = -10439.529

$ws.Range("H65").Value = 11327.513
$ws.Range("I65").Value = 12978.046
$ws.Range("J65").Value = 9191.529
$ws.Range("K65").Value = 64890.23
$ws.Range("L65").Value = 45957.645
$ws.Range("M65").Value = -61770.23
$ws.Range("N65").Value = -52197.645

$ws.Range("H70").Value = 2318.5
$ws.Range("I70").Value = 1194.5
$ws.Range("J70").Value = 2599.5
$ws.Range("K70").Value = 3583.5
$ws.Range("L70").Value = 7798.5
$ws.Range("M70").Value = -3313.5
$ws.Range("N70").Value = -8338.5

$ws.Range("H73").Value = 2318.5
$ws.Range("I73").Value = 1194.5
$ws.Range("J73").Value = 2599.5
$ws.Range("K73").Value = 3583.5
$ws.Range("L73").Value = 7798.5
$ws.Range("M73").Value = -2647.5
$ws.Range("N73").Value = -9670.5

$ws.Range("H98").Value = 1496.5454
$ws.Range("I98").Value = 1496.5454
$ws.Range("K98").Value = 1496.5454
$ws.Range("M98").Value = 1.454600000000028

$ws.Range("H100").Value = 1791.625
$ws.Range("I100").Value = 1266.7693
$ws.Range("K100").Value = 1266.7693
$ws.Range("M100").Value = -725.7692999999999

$ws.Range("H113").Value = 3711
$ws.Range("I113").Value = 3574.7778
$ws.Range("J113").Value = 3805.3076
$ws.Range("K113").Value = 3574.7778
$ws.Range("L113").Value = 3805.3076
$ws.Range("M113").Value = -320.7777999999998
$ws.Range("N113").Value = -10313.3076

$ws.Range("H122").Value = 1496.5454
$ws.Range("I122").Value = 1496.5454
$ws.Range("K122").Value = 4489.6362
$ws.Range("M122").Value = -2039.6362

$ws.Range("H138").Value = 5164.949
$ws.Range("I138").Value = 2862.9048
$ws.Range("J138").Value = 6437.1313
$ws.Range("K138").Value = 8588.714399999999
$ws.Range("L138").Value = 19311.3939
$ws.Range("M138").Value = -3448.714399999999
$ws.Range("N138").Value = -29591.3939

$ws.Range("H141").Value = 6946.75
$ws.Range("I141").Value = 7385.278
$ws.Range("K141").Value = 22155.834
$ws.Range("M141").Value = -16975.834

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H52").Value = 0
$ws.Range("I52").Value = 0
$ws.Range("K52").Value = 0
$ws.Range("M52").ClearContents()

$ws.Range("H97").Value = 13879.9
$ws.Range("I97").Value = 22320
$ws.Range("J97").Value = 5439.8
$ws.Range("K97").Value = 22320
$ws.Range("L97").Value = 5439.8
$ws.Range("M97").Value = -21824
$ws.Range("N97").Value = -6431.8

$ws.Range("H122").Value = 2014.5
$ws.Range("J122").Value = 4500
$ws.Range("L122").Value = 13500
$ws.Range("N122").Value = -18400

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2535.7144
$ws.Range("I86").Value = 2250
$ws.Range("J86").Value = 2750
$ws.Range("K86").Value = 2250
$ws.Range("L86").Value = 2750
$ws.Range("M86").Value = -1127
$ws.Range("N86").Value = -4996

$ws.Range("H89").Value = 2535.7144
$ws.Range("I89").Value = 2250
$ws.Range("J89").Value = 2750
$ws.Range("K89").Value = 11250
$ws.Range("L89").Value = 13750
$ws.Range("M89").Value = -5634
$ws.Range("N89").Value = -24982

$ws.Range("H99").Value = 57594.895
$ws.Range("I99").Value = 95391.17999999999
$ws.Range("K99").Value = 95391.17999999999
$ws.Range("M99").Value = -93893.17999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H60").Value = 50000
$ws.Range("I60").Value = 0
$ws.Range("K60").Value = 0
$ws.Range("M60").ClearContents()

$ws.Range("H94").Value = 3495.8
$ws.Range("I94").Value = 3446.3333
$ws.Range("J94").Value = 3517
$ws.Range("K94").Value = 3446.3333
$ws.Range("L94").Value = 3517
$ws.Range("M94").Value = -2995.3333
$ws.Range("N94").Value = -4419

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 5299.2856
$ws.Range("I55").Value = 5000
$ws.Range("J55").Value = 5349.1665
$ws.Range("K55").Value = 15000
$ws.Range("L55").Value = 16047.4995
$ws.Range("M55").Value = -14823
$ws.Range("N55").Value = -16401.4995

$ws.Range("H132").Value = 1206.5714
$ws.Range("I132").Value = 999.2
$ws.Range("K132").Value = 8992.800000000001
$ws.Range("M132").Value = -6462.800000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3150.2666
$ws.Range("I102").Value = 2250.3076
$ws.Range("K102").Value = 2250.3076
$ws.Range("M102").Value = -628.3076000000001

$ws.Range("H122").Value = 3348.25
$ws.Range("I122").Value = 3312.2856
$ws.Range("K122").Value = 9936.856800000001
$ws.Range("M122").Value = -7486.856800000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2395.9333
$ws.Range("I82").Value = 1155.1
$ws.Range("J82").Value = 3016.35
$ws.Range("K82").Value = 1155.1
$ws.Range("L82").Value = 3016.35
$ws.Range("M82").Value = -794.0999999999999
$ws.Range("N82").Value = -3738.35

$ws.Range("H85").Value = 2395.9333
$ws.Range("I85").Value = 1155.1
$ws.Range("J85").Value = 3016.35
$ws.Range("K85").Value = 1155.1
$ws.Range("L85").Value = 3016.35
$ws.Range("M85").Value = 92.90000000000009
$ws.Range("N85").Value = -5512.35

$ws.Range("H133").Value = 67997
$ws.Range("J133").Value = 67997
$ws.Range("L133").Value = 67997
$ws.Range("N133").Value = -73057

$ws.Range("H138").Value = 110900
$ws.Range("J138").Value = 110900
$ws.Range("L138").Value = 110900
$ws.Range("N138").Value = -121180

$ws.Range("H140").Value = 23929
$ws.Range("J140").Value = 23929
$ws.Range("L140").Value = 23929
$ws.Range("N140").Value = -34289

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 123564.555
$ws.Range("I62").Value = 4801
$ws.Range("J62").Value = 157497
$ws.Range("K62").Value = 4801
$ws.Range("L62").Value = 157497
$ws.Range("M62").Value = -4177
$ws.Range("N62").Value = -158745

$ws.Range("H65").Value = 123564.555
$ws.Range("I65").Value = 4801
$ws.Range("J65").Value = 157497
$ws.Range("K65").Value = 24005
$ws.Range("L65").Value = 787485
$ws.Range("M65").Value = -20885
$ws.Range("N65").Value = -793725

$ws.Range("H133").Value = 88000
$ws.Range("J133").Value = 88000
$ws.Range("L133").Value = 88000
$ws.Range("N133").Value = -98120
